# Apply the crypto price/volume refresh for Sun Feb 18 10:32:20 UTC 2024.
# Column D ("Price") and E ("Volume(1h)") are plain text cells in the sheet,
# and a couple of rows (29/30) swap their Coin/Link too (Cosmos <-> InjectiveProtocol).
#
# Some "Price" values look like plain numbers (e.g. "112.68"), and Excel would
# normally auto-convert those to a numeric cell on assignment. We prefix those
# with a leading apostrophe (the standard Excel 'force text' entry convention)
# so they stay text, matching the original inline-string cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.873.13'
$ws.Range('E2').Value = '  +0.57%  '
$ws.Range('D3').Value = '2.810.98'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').Value = '''356.62'
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').Value = '''112.68'
$ws.Range('E6').Value = '  +3.57%  '
$ws.Range('E7').Value = '  +0.95%  '
$ws.Range('D9').Value = '''0.630'
$ws.Range('E9').Value = '  +8.35%  '
$ws.Range('D10').Value = '''40.40'
$ws.Range('E10').Value = '  +2.58%  '
$ws.Range('E11').Value = '  -0.43%  '
$ws.Range('D12').Value = '''0.0841'
$ws.Range('E12').Value = '  +0.03%  '
$ws.Range('D13').Value = '''19.97'
$ws.Range('E13').Value = '  +3.01%  '
$ws.Range('D14').Value = '''7.79'
$ws.Range('E14').Value = '  +3.21%  '
$ws.Range('D15').Value = '3.256.33'
$ws.Range('E15').Value = '  +1.59%  '
$ws.Range('D16').Value = '2.812.79'
$ws.Range('E16').Value = '  +1.61%  '
$ws.Range('E17').Value = '  +1.70%  '
$ws.Range('D18').Value = '51.871.78'
$ws.Range('E18').Value = '  +0.69%  '
$ws.Range('D19').Value = '''7.64'
$ws.Range('E19').Value = '  +2.96%  '
$ws.Range('D20').Value = '''3.19'
$ws.Range('E20').Value = '  +2.86%  '
$ws.Range('D21').Value = '''13.65'
$ws.Range('D22').Value = '0.0₃0979'
$ws.Range('E22').Value = '  +1.64%  '
$ws.Range('D23').Value = '''70.44'
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('D24').Value = '''268.73'
$ws.Range('E24').Value = '  +0.54%  '
$ws.Range('E25').Value = '  +1.71%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('D27').Value = '''26.17'
$ws.Range('E27').Value = '  -0.41%  '
$ws.Range('E28').Value = '  -0.45%  '
$ws.Range('B29').Value = 'Cosmos'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D29').Value = '''10.42'
$ws.Range('E29').Value = '  +2.69%  '
$ws.Range('B30').Value = 'InjectiveProtocol'
$ws.Range('C30').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = '''38.27'
$ws.Range('E30').Value = '  +11.53%  '
$ws.Range('E31').Value = '  +1.02%  '
$ws.Range('D32').Value = '''6.20'
$ws.Range('E32').Value = '  +0.75%  '
$ws.Range('D33').Value = '''52.45'
$ws.Range('E33').Value = '  +1.85%  '
$ws.Range('D34').Value = '''5.64'
$ws.Range('E34').Value = '  +9.71%  '
$ws.Range('E35').Value = '  -0.05%  '
$ws.Range('D36').Value = '''0.0878'
$ws.Range('E36').Value = '  +5.37%  '
$ws.Range('E37').Value = '  +0.04%  '
$ws.Range('E38').Value = '  -0.10%  '
$ws.Range('E39').Value = '  +3.45%  '
$ws.Range('E40').Value = '  +0.76%  '
$ws.Range('E41').Value = '  +1.54%  '
$ws.Range('D42').Value = '''2.51'
$ws.Range('E42').Value = '  -0.56%  '
$ws.Range('D43').Value = '''120.98'
$ws.Range('E43').Value = '  +1.35%  '
$ws.Range('D44').Value = '''22.04'
$ws.Range('E44').Value = '  +1.79%  '
$ws.Range('E45').Value = '  -0.98%  '
$ws.Range('E46').Value = '  +4.84%  '
$ws.Range('D47').Value = '2.109.06'
$ws.Range('E47').Value = '  +1.34%  '
$ws.Range('D48').Value = '''2.41'
$ws.Range('E48').Value = '  +5.64%  '
$ws.Range('D49').Value = '''0.941'
$ws.Range('E49').Value = '  +0.83%  '
$ws.Range('E50').Value = '  -1.19%  '
$ws.Range('E51').Value = '  +7.81%  '
